$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 857 — everything currently at/after
# row 857 shifts down by one (old 857 -> 858, ..., old 965 -> 966).
$ws.Rows(857).Insert()

# Populate the newly inserted row 857 with its data.
$ws.Range("A857").Value = 3
$ws.Range("B857").Value = "Femacal de La Calera"
$ws.Range("C857").Value = "Coquimbo"
$ws.Range("D857").Value = 45154
$ws.Range("E857").Value = 5
$ws.Range("F857").Value = 100112045
$ws.Range("G857").Value = "Zapallo"
$ws.Range("H857").Value = "Camote"
$ws.Range("I857").Value = "1a (guarda)"
$ws.Range("J857").Value = 165
$ws.Range("K857").Value = 500
$ws.Range("L857").Value = 530
$ws.Range("M857").Value = 512
$ws.Range("N857").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O857").Value = "Región de O'Higgins"
$ws.Range("P857").Value = 512
$ws.Range("Q857").Value = 1
$ws.Range("R857").Value = "Hortaliza"
